# Old script now works with added features:
#  - rename the column header labels (shared strings) used across the
#    shscores*/fi_* sheets from short codes to LaTeX-style labels
#  - update the "mean"/"std" statistics for the "information impact"
#    column on the shscoresUnderwhelming sheet, reflecting the newly
#    computed values from the updated script

$wb = $excel.ActiveWorkbook

$sheetNames = @("shscoresUnderwhelming", "fi_Underwhelming", "shscoresOverwhelming", "fi_Overwhelming")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C1").Value = '$c_i^{deg}$'
    $ws.Range("D1").Value = '$c_i^{betw}$'
    $ws.Range("E1").Value = '$c_i^{ic}$'
    $ws.Range("F1").Value = '$c_i^{ev}$'
}

$wsU = $wb.Worksheets.Item("shscoresUnderwhelming")
$wsU.Range("B2").Value = 0.3666666666666666
$wsU.Range("B3").Value = 0.9637888196533971
